# Auto-generated edit script applying the cryptos.xlsx data refresh
# (commit: "Updated cryptos list on Thu Apr 18 13:57:25 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.444.05'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").Value = '3.035.19'
$ws.Range("E3").Value = '  -0.32%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.71'
$ws.Range("E5").Value = '  +1.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.89'
$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.028.49'
$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("E9").Value = '  +0.54%  '

$ws.Range("E10").Value = '  -0.15%  '

$ws.Range("E11").Value = '  -3.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.446'
$ws.Range("E12").Value = '  -0.54%  '

$ws.Range("E13").Value = '  +0.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.11'
$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("D15").Value = '3.525.49'
$ws.Range("E15").Value = '  -0.65%  '

$ws.Range("D16").Value = '62.568.67'
$ws.Range("E16").Value = '  -0.29%  '

$ws.Range("D17").Value = '3.028.27'
$ws.Range("E17").Value = '  -0.71%  '

$ws.Range("E18").Value = '  -3.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.61'
$ws.Range("E19").Value = '  +0.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '476.30'
$ws.Range("E20").Value = '  +2.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.26'
$ws.Range("E21").Value = '  -0.34%  '

$ws.Range("E22").Value = '  -1.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.05'
$ws.Range("E23").Value = '  +1.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.06'
$ws.Range("E24").Value = '  +3.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.04'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("E26").Value = '  +0.28%  '

$ws.Range("E27").Value = '  +0.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.79'
$ws.Range("E28").Value = '  +0.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.22%  '

$ws.Range("E30").Value = '  +4.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.64'
$ws.Range("E31").Value = '  -0.56%  '

$ws.Range("E32").Value = '  -1.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.37'
$ws.Range("E33").Value = '  +4.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.64'
$ws.Range("E34").Value = '  +4.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '54.89'
$ws.Range("E35").Value = '  -4.57%  '

$ws.Range("E36").Value = '  -0.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '461.96'
$ws.Range("E37").Value = '  +0.45%  '

$ws.Range("D38").Value = '3.155.31'
$ws.Range("E38").Value = '  -1.21%  '

$ws.Range("E39").Value = '  +1.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0388'
$ws.Range("E40").Value = '  -0.44%  '

$ws.Range("E41").Value = '  +1.43%  '

$ws.Range("E42").Value = '  +0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.45'
$ws.Range("E43").Value = '  -1.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.60'
$ws.Range("E44").Value = '  +6.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.245'
$ws.Range("E46").Value = '  -1.17%  '

$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.99'
$ws.Range("E47").Value = '  +1.86%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.108'
$ws.Range("E48").Value = '  +0.08%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.70'
$ws.Range("E49").Value = '  -5.70%  '

$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.0₃0500'
$ws.Range("E50").Value = '  -2.16%  '

$ws.Range("E51").Value = '  +2.28%  '

Write-Output "Applied cryptos list refresh ($([DateTime]::Now))"
